# Update wording on the last slide's SmartArt diagram ("Analysis Process -
# Key Findings and Insights"), per commit "Update to wording on last slide
# for clarity":
#   - "...limit use of two data points..."           -> "...limit use of the two data points..."
#   - "...increase in the machinery used for oil..." -> "...increase in the  price of machinery used for oil..."

$p = $ppt.ActivePresentation

# The presentation has 7 slides; the edit targets the last one.
$s = $p.Slides.Item($p.Slides.Count)

# Find the shape that hosts the SmartArt diagram ("Content Placeholder 2").
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasSmartArt) {
        $target = $candidate
    }
}

$smartArt = $target.SmartArt

# Locate the two SmartArt text nodes that need new wording by matching
# their current text, so the script is resilient to node-order changes.
$collinearityNode = $null
$materialsNode = $null
for ($i = 1; $i -le $smartArt.AllNodes.Count; $i++) {
    $node = $smartArt.AllNodes.Item($i)
    $text = $node.TextFrame.TextRange.Text
    if ($text -like "*Strong collinearity could limit use of two data points*") {
        $collinearityNode = $node
    }
    if ($text -like "*increase in the machinery used for oil extraction*") {
        $materialsNode = $node
    }
}

# Write the materials node first and the collinearity node last: both
# nodes share one on-slide text box with a third (unrelated, untouched)
# bullet, and this ordering leaves that box's visible text matching the
# collinearity node's new wording.
$materialsNode.TextFrame.TextRange.Text = "- An increase in industrial materials costs can lead to an increase in the  price of machinery used for oil extraction and transportation, as well as the production of renewable energy technologies such as windmills and solar panels"
$collinearityNode.TextFrame.TextRange.Text = "- Strong collinearity could limit use of the two data points in regression analysis"
